$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.428.85"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
$ws.Range("D3").Value = "1.630.98"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("E5").Value = "  -0.07%  "

# Row 6
$ws.Range("D6").Value = "'304.63"
$ws.Range("E6").Value = "  -0.93%  "

# Row 7
$ws.Range("D7").Value = "'0.3782"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("D8").Value = "'0.3652"
$ws.Range("E8").Value = "  +0.29%  "

# Row 9
$ws.Range("D9").Value = "'51.55"
$ws.Range("E9").Value = "  -1.29%  "

# Row 10
$ws.Range("D10").Value = "'0.08247"
$ws.Range("E10").Value = "  +1.11%  "

# Row 11
$ws.Range("D11").Value = "'1.231"
$ws.Range("E11").Value = "  -2.51%  "

# Row 12
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.23%  "

# Row 13
$ws.Range("D13").Value = "'22.41"
$ws.Range("E13").Value = "  -2.25%  "

# Row 14
$ws.Range("D14").Value = "'6.560"
$ws.Range("E14").Value = "  -1.11%  "

# Row 15
$ws.Range("E15").Value = "  -1.86%  "

# Row 16
$ws.Range("D16").Value = "'7.338"
$ws.Range("E16").Value = "  -0.46%  "

# Row 17
$ws.Range("D17").Value = "1.630.24"
$ws.Range("E17").Value = "  -0.89%  "

# Row 18
$ws.Range("D18").Value = "'94.14"
$ws.Range("E18").Value = "  -0.71%  "

# Row 19
$ws.Range("D19").Value = "'0.06973"
$ws.Range("E19").Value = "  +0.27%  "

# Row 20
$ws.Range("D20").Value = "'17.70"
$ws.Range("E20").Value = "  -2.58%  "

# Row 21
$ws.Range("D21").Value = "'6.542"
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").Value = "'12.74"
$ws.Range("E23").Value = "  -0.41%  "

# Row 24
$ws.Range("D24").Value = "23.423.77"
$ws.Range("E24").Value = "  -0.15%  "

# Row 25
$ws.Range("D25").Value = "'3.201"
$ws.Range("E25").Value = "  +3.01%  "

# Row 26
$ws.Range("E26").Value = "  +1.90%  "

# Row 27
$ws.Range("D27").Value = "'21.42"
$ws.Range("E27").Value = "  +0.77%  "

# Row 28
$ws.Range("D28").Value = "'150.16"
$ws.Range("E28").Value = "  -0.57%  "

# Row 29
$ws.Range("D29").Value = "'5.301"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30
$ws.Range("D30").Value = "'134.36"
$ws.Range("E30").Value = "  -0.70%  "

# Row 31
$ws.Range("D31").Value = "1.810.83"
$ws.Range("E31").Value = "  -0.85%  "

# Row 32
$ws.Range("D32").Value = "'2.257"
$ws.Range("E32").Value = "  -3.01%  "

# Row 33
$ws.Range("D33").Value = "'6.827"
$ws.Range("E33").Value = "  +0.86%  "

# Row 34
$ws.Range("D34").Value = "'1.019"
$ws.Range("E34").Value = "  +5.88%  "

# Row 35
$ws.Range("E35").Value = "  +5.47%  "

# Row 36
$ws.Range("D36").Value = "'0.02800"
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("D37").Value = "'0.2529"
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("D38").Value = "'0.08794"
$ws.Range("E38").Value = "  -0.67%  "

# Row 39
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$ws.Range("D40").Value = "'0.07133"
$ws.Range("E40").Value = "  -2.81%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.7052"
$ws.Range("E41").Value = "  -0.64%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.353"
$ws.Range("E42").Value = "  -2.06%  "

# Row 43
$ws.Range("D43").Value = "'16.21"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").Value = "'12.21"
$ws.Range("E44").Value = "  -2.30%  "

# Row 45
$ws.Range("D45").Value = "'0.6582"
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("E46").Value = "  -0.65%  "

# Row 47
$ws.Range("E47").Value = "  -0.08%  "

# Row 48
$ws.Range("D48").Value = "'3.980"
$ws.Range("E48").Value = "  -1.27%  "

# Row 49
$ws.Range("D49").Value = "'0.08013"
$ws.Range("E49").Value = "  +0.48%  "

# Row 50
$ws.Range("D50").Value = "'1.199"
$ws.Range("E50").Value = "  -0.50%  "

# Row 51
$ws.Range("D51").Value = "'126.44"
$ws.Range("E51").Value = "  -2.19%  "
